# add functionality for rsm
# Refresh the CNH call-code table: row labels shift (new codes CNH11/CNH15/CNH46
# inserted, CNH24/CNH26 retired) and all call counts (columns B:J) are updated,
# including a brand-new row 32 (CNH66).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @("CNH11", 0, 0, 0, 0, 0, 0, 0, 0, 0)
  3 = @("CNH12", 0, 0, 6, 0, 0, 0, 1, 1, 5)
  4 = @("CNH13", 3, 6, 5, 2, 5, 1, 2, 2, 5)
  5 = @("CNH14", 0, 1, 0, 0, 0, 1, 0, 0, 0)
  6 = @("CNH15", 0, 0, 0, 0, 2, 0, 5, 0, 0)
  7 = @("CNH16", 4, 3, 7, 1, 0, 6, 3, 3, 1)
  8 = @("CNH17", 3, 3, 6, 6, 5, 7, 3, 0, 3)
  9 = @("CNH21", 14, 10, 6, 13, 17, 13, 7, 12, 15)
  10 = @("CNH22", 3, 2, 2, 0, 0, 3, 2, 0, 2)
  11 = @("CNH23", 3, 2, 0, 1, 0, 0, 1, 1, 0)
  12 = @("CNH27", 0, 0, 0, 0, 0, 0, 0, 0, 0)
  13 = @("CNH31", 2, 1, 2, 1, 3, 6, 6, 1, 0)
  14 = @("CNH32", 5, 2, 6, 3, 4, 4, 5, 0, 5)
  15 = @("CNH33", 1, 1, 5, 0, 2, 5, 4, 1, 0)
  16 = @("CNH34", 4, 1, 3, 1, 9, 16, 10, 4, 0)
  17 = @("CNH35", 4, 2, 3, 0, 0, 3, 4, 0, 3)
  18 = @("CNH41", 9, 6, 6, 5, 2, 3, 7, 2, 2)
  19 = @("CNH42", 5, 6, 3, 6, 5, 4, 12, 2, 2)
  20 = @("CNH44", 4, 2, 1, 1, 0, 7, 8, 0, 1)
  21 = @("CNH45", 13, 6, 5, 9, 12, 8, 9, 9, 4)
  22 = @("CNH46", 0, 1, 1, 0, 0, 1, 0, 0, 1)
  23 = @("CNH51", 2, 1, 5, 2, 5, 8, 4, 0, 3)
  24 = @("CNH53", 9, 9, 2, 4, 7, 4, 1, 0, 0)
  25 = @("CNH54", 3, 2, 4, 1, 1, 6, 10, 0, 10)
  26 = @("CNH55", 6, 3, 6, 0, 0, 2, 2, 0, 2)
  27 = @("CNH56", 11, 3, 13, 4, 8, 14, 9, 2, 4)
  28 = @("CNH57", 11, 19, 10, 11, 22, 20, 14, 7, 20)
  29 = @("CNH63", 7, 1, 7, 0, 3, 6, 7, 6, 8)
  30 = @("CNH64", 4, 11, 1, 4, 6, 0, 10, 11, 1)
  31 = @("CNH65", 11, 1, 1, 4, 6, 3, 5, 3, 4)
  32 = @("CNH66", 7, 2, 8, 6, 8, 15, 11, 3, 2)
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  for ($c = 1; $c -le $row.Length; $c++) {
    $ws.Cells.Item([int]$r, $c).Value = $row[$c - 1]
  }
}
